$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 41, shifting existing rows 41-54 down to 42-55.
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new weekly record.
$ws.Range("A41").Value2 = 5
$ws.Range("B41").Value2 = "Macroferia Regional de Talca"
$ws.Range("C41").Value2 = "Maule"
$ws.Range("D41").Value2 = 44932
$ws.Range("E41").Value2 = 7
$ws.Range("F41").Value2 = "Fruta"
$ws.Range("G41").Value2 = 100103
$ws.Range("H41").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I41").Value2 = 100103003
$ws.Range("J41").Value2 = "Damasco"
$ws.Range("K41").Value2 = "Dina"
$ws.Range("L41").Value2 = "Primera"
$ws.Range("M41").Value2 = 150
$ws.Range("N41").Value2 = 10000
$ws.Range("O41").Value2 = 10000
$ws.Range("P41").Value2 = 10000
$ws.Range("Q41").Value2 = "$/bandeja 10 kilos"
$ws.Range("R41").Value2 = "Región de O'Higgins"
$ws.Range("S41").Value2 = 1000
$ws.Range("T41").Value2 = 10
